$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name
$ws.Name = "Through 2021-12-21"

# Update the December label (row 13, column A)
$ws.Range("A13").Value = "December (through 12-21)"

# Update December row (row 13) values
$ws.Range("B13").Value = 28
$ws.Range("C13").Value = 66
$ws.Range("D13").Value = 81
$ws.Range("E13").Value = 47
$ws.Range("F13").Value = 39
$ws.Range("G13").Value = 98
$ws.Range("H13").Value = 139

# Update Total row (row 14) values
$ws.Range("B14").Value = 319
$ws.Range("C14").Value = 629
$ws.Range("D14").Value = 902
$ws.Range("E14").Value = 729
$ws.Range("F14").Value = 573
$ws.Range("G14").Value = 1362
$ws.Range("H14").Value = 1782
